$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"

$ws.Range('D2').Value = $q + '65.099.27'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = $q + '  +1.64%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = $q + '3.184.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = $q + '  +4.07%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = $q + '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = $q + '  -0.08%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = $q + '577.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = $q + '  +2.92%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = $q + '151.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = $q + '  +5.50%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = $q + '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = $q + '  -0.04%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = $q + '3.184.06'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = $q + '  +4.08%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = $q + '0.530'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = $q + '  +3.13%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = $q + '0.164'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = $q + '  +5.61%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = $q + '6.23'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = $q + '  +2.03%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = $q + '0.508'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = $q + '  +2.90%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = $q + '  +18.38%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = $q + '38.15'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = $q + '  +7.46%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = $q + '3.704.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = $q + '  +4.05%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = $q + '65.158.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = $q + '  +1.68%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = $q + '3.180.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = $q + '  +3.88%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = $q + '  +6.33%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = $q + '  +1.04%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = $q + '514.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = $q + '  +7.39%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = $q + '14.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = $q + '  +7.14%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = $q + '0.737'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = $q + '  +8.07%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = $q + '15.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = $q + '  +6.88%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = $q + '7.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = $q + '  +3.53%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = $q + '85.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = $q + '  +3.66%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = $q + '  +0.10%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = $q + '  +12.00%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = $q + '2.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = $q + '  +4.02%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = $q + '  +8.11%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = $q + '  +7.16%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = $q + '  +13.82%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = $q + 'Mantle'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = $q + 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = $q + '1.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = $q + '  +6.59%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = $q + 'FirstDigitalUSD'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = $q + 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = $q + '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = $q + '  -0.13%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = $q + '6.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = $q + '  +9.01%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = $q + '6.67'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = $q + '  +7.07%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = $q + '55.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = $q + '  +1.52%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = $q + '0.0900'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = $q + '  +10.65%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = $q + '478.74'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = $q + '  +7.19%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = $q + '3.15'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = $q + '  +10.94%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = $q + '  +2.49%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = $q + 'Cosmos'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = $q + 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = $q + '8.68'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = $q + '  +4.79%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = $q + 'Maker'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = $q + 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = $q + '3.073.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = $q + '  +1.80%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = $q + '  +3.77%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = $q + '0.289'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = $q + '  +9.46%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = $q + '  +11.24%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = $q + '29.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = $q + '  +5.55%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = $q + '0.0₃0604'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = $q + '  +16.90%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E49').Value = $q + '  +2.31%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = $q + '  +10.40%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = $q + '121.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = $q + '  +2.52%  '
$ws.Range('E51').Style = 'Normal'

Write-Output "done"
